$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22: LeetCode 142 - 单链表中的环2
$ws.Range("A22").Value = 142
$ws.Range("B22").Value = "单链表中的环2"
$ws.Range("D22").Value = "||"
$ws.Range("E22").Value = "同上一题差不多"

# New row 23: LeetCode 143 - 重排链表
$ws.Range("A23").Value = 143
$ws.Range("B23").Value = "重排链表"
$ws.Range("D23").Value = "||"
$ws.Range("E23").Value = "不是很难"

# New row 24: LeetCode 144 - 树的前序遍历 (no note column)
$ws.Range("A24").Value = 144
$ws.Range("B24").Value = "树的前序遍历"
$ws.Range("D24").Value = "|"

# New row 25: LeetCode 145 - 树的后续遍历
$ws.Range("A25").Value = 145
$ws.Range("B25").Value = "树的后续遍历"
$ws.Range("D25").Value = "||||"
$ws.Range("E25").Value = "迭代法一定要会！！！"

# Move the active selection to match the author's final cursor position
$ws.Range("E25").Select()
